# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Update the "data" sheet's time_taken (F) column timestamps ---
$ws1.Range("F2").Value = "2021-10-05 14:35:46.730439"
$ws1.Range("F3").Value = "2021-10-05 14:35:46.730447"
$ws1.Range("F4").Value = "2021-10-05 14:35:46.730450"
$ws1.Range("F5").Value = "2021-10-05 14:35:46.730453"
$ws1.Range("F6").Value = "2021-10-05 14:35:46.730456"
$ws1.Range("F7").Value = "2021-10-05 14:35:46.730458"

# --- 2. Add a new "metadata" worksheet after "data" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# --- 3. Copy the bold/border/centered header style from sheet1!B1 onto the
#        header row and the A2 "index" cell, without creating new style
#        entries (PasteSpecial Formats reuses the existing style record). ---
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Header row values ---
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# --- 5. Data row values ---
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Spondylocostal Dysostosis"
$ws2.Range("C2").Value = 177

# D2 must be stored as the literal text "0.5" (not a number). Writing the
# string directly via .Value lets Excel's numeric auto-detection convert it
# to a number, so instead build it as a text formula in a scratch cell and
# paste-special just the (text-typed) value into D2.
$scratch = $ws2.Range("Z100")
$scratch.Formula = '="0.5"'
$scratch.Copy()
$ws2.Range("D2").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false

$ws2.Range("E2").Value = "2020-11-13T07:48:12.323235Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:46.726640"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/177/?format=json"

# --- 6. Leave the original "data" sheet as the active/selected tab, since
#        the diff leaves <bookViews>/activeTab untouched. ---
$ws1.Activate() | Out-Null
